$d = $word.ActiveDocument

# Locate the "Assignment 2.3" heading text so we don't depend on hardcoded
# character offsets. After Find.Execute the search range collapses to the
# matched text, so $headerRange.End is the position right after the "3".
$headerRange = $d.Content
$findResult = $headerRange.Find.Execute("Assignment 2.3", $true, $false, $false,
                                         $false, $false, $true, 1, $false, "", 0)
$afterDigit = $headerRange.End
$digitStart = $afterDigit - 1
$digitRange = $d.Range($digitStart, $afterDigit)

# This interop engine coalesces adjacent runs that share identical
# formatting whenever a paragraph is edited - which here would merge the
# digit run with the "." run before it and the " - One-Way Binding" runs
# after it, flattening the run structure the target keeps split apart.
# Bookmarks act as hard boundaries that block that coalescing, so we
# bracket the edit with them before changing the text.

# Temporary boundary right before the digit, so it doesn't get folded back
# into the preceding "." run.
$beforeDigit = $d.Range($digitStart, $digitStart)
$d.Bookmarks.Add("ztmp_boundary", $beforeDigit)

# The document's "_GoBack" bookmark currently brackets the image paragraph
# further down; move it to sit right after the assignment number instead
# (between the digit and " - One-Way Binding"). Adding a bookmark under a
# name that already exists relocates it, so this both removes it from the
# image paragraph and stops the trailing runs from merging into the digit.
$goBackRange = $d.Range($afterDigit, $afterDigit)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Now make the actual content change: "3" -> "2" (Assignment 2.3 -> 2.2).
$digitRange.Text = "2"

# Drop the temporary boundary bookmark - the runs are already split apart
# by now, so removing it will not cause anything to re-merge.
$tmp = $d.Bookmarks("ztmp_boundary")
$tmp.Delete()
